$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$updates = @(
    @{Addr="J6"; Value=499},
    @{Addr="L6"; Value=1497},
    @{Addr="I6"; Value=7983.4546},
    @{Addr="M6"; Value=-23838.3638},
    @{Addr="K6"; Value=23950.3638},
    @{Addr="N6"; Value=-1721},
    @{Addr="H6"; Value=7359.75},
    @{Addr="J19"; Value=3248.5},
    @{Addr="N19"; Value=-3598.5},
    @{Addr="H19"; Value=2290.5},
    @{Addr="I19"; Value=374.5},
    @{Addr="M19"; Value=-199.5},
    @{Addr="K19"; Value=374.5},
    @{Addr="L19"; Value=3248.5},
    @{Addr="H69"; Value=19163.5},
    @{Addr="J69"; Value=21197.4},
    @{Addr="N69"; Value=-65340.2},
    @{Addr="L69"; Value=63592.2},
    @{Addr="J72"; Value=21197.4},
    @{Addr="N72"; Value=-199512.6},
    @{Addr="L72"; Value=190776.6},
    @{Addr="H72"; Value=19163.5},
    @{Addr="L100"; Value=25437.5},
    @{Addr="M100"; Value=-740.8462},
    @{Addr="H100"; Value=4502.6},
    @{Addr="K100"; Value=1281.8462},
    @{Addr="J100"; Value=25437.5},
    @{Addr="N100"; Value=-26519.5},
    @{Addr="I100"; Value=1281.8462},
    @{Addr="N112"; Value=-8888.6921},
    @{Addr="H112"; Value=2143.8572},
    @{Addr="J112"; Value=2224.2307},
    @{Addr="L112"; Value=6672.6921},
    @{Addr="K132"; Value=2326.32348},
    @{Addr="H132"; Value=1329.25},
    @{Addr="M132"; Value=203.67652},
    @{Addr="I132"; Value=775.44116},
    @{Addr="K135"; Value=12236.9994},
    @{Addr="H135"; Value=1479.7142},
    @{Addr="M135"; Value=-9701.999400000001},
    @{Addr="I135"; Value=1359.6666},
    @{Addr="H137"; Value=6066.9},
    @{Addr="I137"; Value=5500.3213},
    @{Addr="K137"; Value=16500.9639},
    @{Addr="M137"; Value=-13950.9639}
)
foreach ($u in $updates) {
    if ($u.Clear) {
        $ws.Range($u.Addr).ClearContents()
    } else {
        $ws.Range($u.Addr).Value = $u.Value
    }
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$updates = @(
    @{Addr="L7"; Value=61598},
    @{Addr="H7"; Value=66855.71000000001},
    @{Addr="N7"; Value=-61826},
    @{Addr="J7"; Value=61598},
    @{Addr="H61"; Value=34097180},
    @{Addr="I61"; Value=33338064},
    @{Addr="K61"; Value=33338064},
    @{Addr="M61"; Value=-33337852},
    @{Addr="I74"; Value=13890176},
    @{Addr="K74"; Value=13890176},
    @{Addr="M74"; Value=-13889302},
    @{Addr="H74"; Value=11305779},
    @{Addr="M77"; Value=-69446512},
    @{Addr="I77"; Value=13890176},
    @{Addr="K77"; Value=69450880},
    @{Addr="H77"; Value=11305779},
    @{Addr="J92"; Value=62341.145},
    @{Addr="H92"; Value=62341.145},
    @{Addr="N92"; Value=-67333.14499999999},
    @{Addr="L92"; Value=62341.145},
    @{Addr="J121"; Value=60127.5},
    @{Addr="H121"; Value=60127.5},
    @{Addr="N121"; Value=-63621.5},
    @{Addr="L121"; Value=60127.5},
    @{Addr="K136"; Value=100014192},
    @{Addr="M136"; Value=-100011642},
    @{Addr="I136"; Value=33338064},
    @{Addr="H136"; Value=34097180},
    @{Addr="N138"; Value=-410279},
    @{Addr="J138"; Value=399999},
    @{Addr="L138"; Value=399999},
    @{Addr="H138"; Value=399999}
)
foreach ($u in $updates) {
    if ($u.Clear) {
        $ws.Range($u.Addr).ClearContents()
    } else {
        $ws.Range($u.Addr).Value = $u.Value
    }
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$updates = @(
    @{Addr="K86"; Value=2799.2354},
    @{Addr="M86"; Value=-1676.2354},
    @{Addr="H86"; Value=2398},
    @{Addr="I86"; Value=2799.2354},
    @{Addr="K89"; Value=13996.177},
    @{Addr="H89"; Value=2398},
    @{Addr="M89"; Value=-8380.177},
    @{Addr="I89"; Value=2799.2354},
    @{Addr="M134"; Value=-2683.3125},
    @{Addr="H134"; Value=613046.9399999999},
    @{Addr="K134"; Value=5218.3125},
    @{Addr="I134"; Value=1739.4375}
)
foreach ($u in $updates) {
    if ($u.Clear) {
        $ws.Range($u.Addr).ClearContents()
    } else {
        $ws.Range($u.Addr).Value = $u.Value
    }
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$updates = @(
    @{Addr="I31"; Value=2317.5715},
    @{Addr="J31"; Value=2916308.5},
    @{Addr="K31"; Value=2317.5715},
    @{Addr="N31"; Value=-2916898.5},
    @{Addr="H31"; Value=806177.1},
    @{Addr="L31"; Value=2916308.5},
    @{Addr="M31"; Value=-2022.5715},
    @{Addr="I34"; Value=2317.5715},
    @{Addr="K34"; Value=2317.5715},
    @{Addr="N34"; Value=-2916712.5},
    @{Addr="L34"; Value=2916308.5},
    @{Addr="M34"; Value=-2115.5715},
    @{Addr="H34"; Value=806177.1},
    @{Addr="J34"; Value=2916308.5},
    @{Addr="I58"; Value=3264.5},
    @{Addr="N58"; Value=-8910.666999999999},
    @{Addr="J58"; Value=8504.666999999999},
    @{Addr="K58"; Value=3264.5},
    @{Addr="L58"; Value=8504.666999999999},
    @{Addr="H58"; Value=4693.636},
    @{Addr="M58"; Value=-3061.5},
    @{Addr="J92"; Value=68495},
    @{Addr="H92"; Value=68495},
    @{Addr="N92"; Value=-73487},
    @{Addr="L92"; Value=68495},
    @{Addr="M99"; Value=-1350.6667},
    @{Addr="H99"; Value=3074.75},
    @{Addr="K99"; Value=2848.6667},
    @{Addr="I99"; Value=2848.6667},
    @{Addr="M126"; Value=-6076.000100000001},
    @{Addr="I126"; Value=2848.6667},
    @{Addr="K126"; Value=8546.000100000001},
    @{Addr="H126"; Value=3074.75},
    @{Addr="K136"; Value=9793.5},
    @{Addr="M136"; Value=-7243.5},
    @{Addr="J136"; Value=8504.666999999999},
    @{Addr="I136"; Value=3264.5},
    @{Addr="L136"; Value=25514.001},
    @{Addr="H136"; Value=4693.636},
    @{Addr="N136"; Value=-30614.001}
)
foreach ($u in $updates) {
    if ($u.Clear) {
        $ws.Range($u.Addr).ClearContents()
    } else {
        $ws.Range($u.Addr).Value = $u.Value
    }
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$updates = @(
    @{Addr="I4"; Value=22510660},
    @{Addr="K4"; Value=67531980},
    @{Addr="H4"; Value=43289936},
    @{Addr="M4"; Value=-67531868},
    @{Addr="H7"; Value=1120.8},
    @{Addr="K7"; Value=1201.5},
    @{Addr="M7"; Value=-1089.5},
    @{Addr="I7"; Value=400.5},
    @{Addr="I50"; Value=97.666664},
    @{Addr="L50"; Value=1311.5625},
    @{Addr="H50"; Value=383.57895},
    @{Addr="K50"; Value=292.999992},
    @{Addr="M50"; Value=188.000008},
    @{Addr="N50"; Value=-2273.5625},
    @{Addr="J50"; Value=437.1875},
    @{Addr="J53"; Value=437.1875},
    @{Addr="K53"; Value=292.999992},
    @{Addr="M53"; Value=188.000008},
    @{Addr="H53"; Value=383.57895},
    @{Addr="I53"; Value=97.666664},
    @{Addr="N53"; Value=-2273.5625},
    @{Addr="L53"; Value=1311.5625},
    @{Addr="H92"; Value=581.6667},
    @{Addr="M92"; Value=-469.5},
    @{Addr="I92"; Value=572.5},
    @{Addr="K92"; Value=1717.5},
    @{Addr="M102"; Clear=$true},
    @{Addr="K102"; Value=0},
    @{Addr="H102"; Value=6000},
    @{Addr="I102"; Value=0},
    @{Addr="H107"; Value=832.26666},
    @{Addr="N107"; Value=-6336.79998},
    @{Addr="J107"; Value=832.26666},
    @{Addr="L107"; Value=2496.79998},
    @{Addr="N131"; Value=-37521.144},
    @{Addr="J131"; Value=9147.048000000001},
    @{Addr="H131"; Value=119746.02},
    @{Addr="L131"; Value=27441.144},
    @{Addr="M139"; Value=-11539.5005},
    @{Addr="K139"; Value=16679.5005},
    @{Addr="I139"; Value=5559.8335},
    @{Addr="H139"; Value=5279.9165}
)
foreach ($u in $updates) {
    if ($u.Clear) {
        $ws.Range($u.Addr).ClearContents()
    } else {
        $ws.Range($u.Addr).Value = $u.Value
    }
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$updates = @(
    @{Addr="K52"; Value=0},
    @{Addr="H52"; Value=0},
    @{Addr="M52"; Clear=$true},
    @{Addr="I52"; Value=0},
    @{Addr="K132"; Value=600002160},
    @{Addr="H132"; Value=166683940},
    @{Addr="M132"; Value=-599999630},
    @{Addr="I132"; Value=200000720},
    @{Addr="N141"; Value=-50360},
    @{Addr="L141"; Value=40000},
    @{Addr="H141"; Value=40000},
    @{Addr="J141"; Value=40000}
)
foreach ($u in $updates) {
    if ($u.Clear) {
        $ws.Range($u.Addr).ClearContents()
    } else {
        $ws.Range($u.Addr).Value = $u.Value
    }
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$updates = @(
    @{Addr="K40"; Value=3638.0667},
    @{Addr="M40"; Value=-3502.0667},
    @{Addr="I40"; Value=3638.0667},
    @{Addr="H40"; Value=4252.346},
    @{Addr="N61"; Clear=$true},
    @{Addr="J61"; Value=0},
    @{Addr="H61"; Value=3272.3333},
    @{Addr="I61"; Value=3272.3333},
    @{Addr="K61"; Value=3272.3333},
    @{Addr="M61"; Value=-3070.3333},
    @{Addr="L61"; Value=0},
    @{Addr="K113"; Value=3272.3333},
    @{Addr="I113"; Value=3272.3333},
    @{Addr="J113"; Value=0},
    @{Addr="M113"; Value=-1102.3333},
    @{Addr="N113"; Clear=$true},
    @{Addr="H113"; Value=3272.3333},
    @{Addr="L113"; Value=0},
    @{Addr="K136"; Value=53497.99800000001},
    @{Addr="M136"; Value=-50947.99800000001},
    @{Addr="I136"; Value=17832.666},
    @{Addr="H136"; Value=99718.69}
)
foreach ($u in $updates) {
    if ($u.Clear) {
        $ws.Range($u.Addr).ClearContents()
    } else {
        $ws.Range($u.Addr).Value = $u.Value
    }
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$updates = @(
    @{Addr="J80"; Value=46483.332},
    @{Addr="N80"; Value=-48479.332},
    @{Addr="L80"; Value=46483.332},
    @{Addr="H80"; Value=39886},
    @{Addr="N83"; Value=-149433.996},
    @{Addr="L83"; Value=139449.996},
    @{Addr="J83"; Value=46483.332},
    @{Addr="H83"; Value=39886},
    @{Addr="L100"; Value=1618.5},
    @{Addr="M100"; Value=-786.3334},
    @{Addr="H100"; Value=746.8570999999999},
    @{Addr="K100"; Value=1327.3334},
    @{Addr="J100"; Value=809.25},
    @{Addr="N100"; Value=-2700.5},
    @{Addr="I100"; Value=663.6667},
    @{Addr="N132"; Value=-3765282.5},
    @{Addr="L132"; Value=3760222.5},
    @{Addr="K132"; Value=11878.125},
    @{Addr="H132"; Value=420442.1},
    @{Addr="M132"; Value=-9348.125},
    @{Addr="I132"; Value=3959.375},
    @{Addr="J132"; Value=1253407.5},
    @{Addr="K136"; Value=17906.4552},
    @{Addr="M136"; Value=-15356.4552},
    @{Addr="I136"; Value=5968.8184},
    @{Addr="H136"; Value=6039.1875},
    @{Addr="H140"; Value=42721},
    @{Addr="N140"; Value=-53081},
    @{Addr="J140"; Value=42721},
    @{Addr="L140"; Value=42721}
)
foreach ($u in $updates) {
    if ($u.Clear) {
        $ws.Range($u.Addr).ClearContents()
    } else {
        $ws.Range($u.Addr).Value = $u.Value
    }
}
